$wb = $excel.ActiveWorkbook

# --- About sheet: clear the stray date value in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Clear()

# --- RTMF-passengers sheet: update mode-shift fractions for LDVs row ---
$wsPass = $wb.Worksheets.Item("RTMF-passengers")
$wsPass.Range("C2").Value = 0.33
$wsPass.Range("E2").Value = 0.33
$wsPass.Range("I2").Formula = "=1-SUM(B2:G2)"
